$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh inserts a new pair of rows (Zafiro rojo / Zafiro verde)
# at the top of this variety block, pushing the existing rows 312-329 down
# to 314-331.
$ws.Rows("312:313").Insert()

# Row 312 - new "Zafiro rojo" entry for the new week (date 44826)
$ws.Cells.Item(312, 1).Value = 7
$ws.Cells.Item(312, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(312, 3).Value = "Ñuble"
$ws.Cells.Item(312, 4).Value = 44826
$ws.Cells.Item(312, 5).Value = 16
$ws.Cells.Item(312, 6).Value = 100112002
$ws.Cells.Item(312, 7).Value = "Pimiento"
$ws.Cells.Item(312, 8).Value = "Zafiro rojo"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 120
$ws.Cells.Item(312, 11).Value = 16000
$ws.Cells.Item(312, 12).Value = 17000
$ws.Cells.Item(312, 13).Value = 16500
$ws.Cells.Item(312, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(312, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(312, 16).Value = 1100
$ws.Cells.Item(312, 17).Value = 15
$ws.Cells.Item(312, 18).Value = "Hortaliza"

# Row 313 - new "Zafiro verde" entry for the new week (date 44826)
$ws.Cells.Item(313, 1).Value = 7
$ws.Cells.Item(313, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(313, 3).Value = "Ñuble"
$ws.Cells.Item(313, 4).Value = 44826
$ws.Cells.Item(313, 5).Value = 16
$ws.Cells.Item(313, 6).Value = 100112002
$ws.Cells.Item(313, 7).Value = "Pimiento"
$ws.Cells.Item(313, 8).Value = "Zafiro verde"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 120
$ws.Cells.Item(313, 11).Value = 16000
$ws.Cells.Item(313, 12).Value = 17000
$ws.Cells.Item(313, 13).Value = 16500
$ws.Cells.Item(313, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(313, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(313, 16).Value = 1100
$ws.Cells.Item(313, 17).Value = 15
$ws.Cells.Item(313, 18).Value = "Hortaliza"
